$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '31.024.67'
$ws.Range("E2").Value = '  +1.25%  '

# Row 3
$ws.Range("D3").Value = '1.955.26'
$ws.Range("E3").Value = '  -0.15%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''245.15'
$ws.Range("E5").Value = '  -1.23%  '

# Row 6
$ws.Range("E6").Value = '  +0.08%  '

# Row 7
$ws.Range("D7").Value = '''0.4858'
$ws.Range("E7").Value = '  +0.72%  '

# Row 8
$ws.Range("E8").Value = '  +0.26%  '

# Row 9
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").Value = '''19.17'
$ws.Range("E10").Value = '  -1.14%  '

# Row 11
$ws.Range("D11").Value = '''107.21'
$ws.Range("E11").Value = '  -3.15%  '

# Row 12
$ws.Range("D12").Value = '1.951.64'
$ws.Range("E12").Value = '  -0.41%  '

# Row 13
$ws.Range("D13").Value = '''0.07813'
$ws.Range("E13").Value = '  +1.15%  '

# Row 14
$ws.Range("D14").Value = '''5.454'
$ws.Range("E14").Value = '  -0.10%  '

# Row 15
$ws.Range("D15").Value = '''0.7021'
$ws.Range("E15").Value = '  +2.15%  '

# Row 16
$ws.Range("D16").Value = '''283.32'
$ws.Range("E16").Value = '  -2.96%  '

# Row 17
$ws.Range("D17").Value = '31.045.87'
$ws.Range("E17").Value = '  +1.25%  '

# Row 18
$ws.Range("D18").Value = '''13.18'
$ws.Range("E18").Value = '  -0.62%  '

# Row 19
$ws.Range("D19").Value = '''0.000007684'
$ws.Range("E19").Value = '  +0.14%  '

# Row 20
$ws.Range("D20").Value = '2.209.29'
$ws.Range("E20").Value = '  -0.48%  '

# Row 21
$ws.Range("E21").Value = '  +0.16%  '

# Row 22
$ws.Range("D22").Value = '''5.491'
$ws.Range("E22").Value = '  -2.75%  '

# Row 23
$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").Value = '''6.492'
$ws.Range("E24").Value = '  -1.60%  '

# Row 25
$ws.Range("D25").Value = '''9.805'
$ws.Range("E25").Value = '  -0.92%  '

# Row 26
$ws.Range("D26").Value = '''169.93'
$ws.Range("E26").Value = '  -0.74%  '

# Row 27
$ws.Range("D27").Value = '''19.94'
$ws.Range("E27").Value = '  -0.90%  '

# Row 28
$ws.Range("D28").Value = '''2.197'
$ws.Range("E28").Value = '  +0.28%  '

# Row 29
$ws.Range("E29").Value = '  -1.78%  '

# Row 30
$ws.Range("D30").Value = '''1.407'
$ws.Range("E30").Value = '  -2.16%  '

# Row 31
$ws.Range("D31").Value = '''1.582'
$ws.Range("E31").Value = '  -1.12%  '

# Row 32
$ws.Range("D32").Value = '''4.609'
$ws.Range("E32").Value = '  -1.87%  '

# Row 33
$ws.Range("D33").Value = '''4.436'
$ws.Range("E33").Value = '  -0.03%  '

# Row 34
$ws.Range("D34").Value = '''0.04921'
$ws.Range("E34").Value = '  -3.75%  '

# Row 35
$ws.Range("D35").Value = '''0.7628'
$ws.Range("E35").Value = '  -1.90%  '

# Row 36
$ws.Range("D36").Value = '''1.171'
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("D37").Value = '''2.727'
$ws.Range("E37").Value = '  -0.23%  '

# Row 38
$ws.Range("D38").Value = '''0.02004'
$ws.Range("E38").Value = '  -2.61%  '

# Row 39
$ws.Range("D39").Value = '''2.705'
$ws.Range("E39").Value = '  -0.19%  '

# Row 40
$ws.Range("D40").Value = '''6.535'
$ws.Range("E40").Value = '  +5.86%  '

# Row 41
$ws.Range("D41").Value = '''2.098'
$ws.Range("E41").Value = '  +1.32%  '

# Row 42
$ws.Range("D42").Value = '''74.93'
$ws.Range("E42").Value = '  +7.06%  '

# Row 43
$ws.Range("D43").Value = '''0.8875'
$ws.Range("E43").Value = '  +1.57%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''109.24'
$ws.Range("E44").Value = '  -0.66%  '

# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.4456'
$ws.Range("E45").Value = '  -0.14%  '

# Row 46
$ws.Range("D46").Value = '''8.181'
$ws.Range("E46").Value = '  +10.69%  '

# Row 47
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("D48").Value = '1.002.01'
$ws.Range("E48").Value = '  +10.70%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''9.339'
$ws.Range("E49").Value = '  -0.86%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.1255'
$ws.Range("E50").Value = '  -1.89%  '

# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '''35.72'
$ws.Range("E51").Value = '  -0.57%  '
